$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "plain" string (never mis-parsed as a number,
# e.g. price strings using dots as thousands separators, percent strings
# with surrounding spaces, coin names, or URLs) -- safe to assign directly.
$plainUpdates = @{
    "D2"  = "66.232.67"
    "E2"  = "  -1.25%  "
    "D3"  = "3.773.57"
    "E3"  = "  +1.70%  "
    "E4"  = "  -0.27%  "
    "E5"  = "  -3.10%  "
    "E6"  = "  +1.23%  "
    "D7"  = "3.761.51"
    "E7"  = "  +1.56%  "
    "E8"  = "  -5.54%  "
    "E9"  = "  +0.08%  "
    "E10" = "  -4.81%  "
    "E11" = "  -8.85%  "
    "E12" = "  -9.75%  "
    "E13" = "  -4.66%  "
    "D14" = "4.357.79"
    "E14" = "  +1.13%  "
    "E15" = "  -4.27%  "
    "E16" = "  +13.25%  "
    "E17" = "  -1.49%  "
    "D18" = "3.785.45"
    "E18" = "  +1.93%  "
    "E19" = "  -5.53%  "
    "D20" = "66.435.89"
    "E20" = "  -1.08%  "
    "E21" = "  -5.87%  "
    "E22" = "  -8.22%  "
    "E23" = "  -9.77%  "
    "E24" = "  -5.09%  "
    "E25" = "  -3.17%  "
    "E26" = "  +14.65%  "
    "E27" = "  -5.91%  "
    "E28" = "  -6.06%  "
    "E29" = "  -9.20%  "
    "E30" = "  +13.05%  "
    "E31" = "  -1.84%  "
    "E32" = "  -1.54%  "
    "B33" = "RenderToken"
    "C33" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "E33" = "  +1.12%  "
    "B34" = "Toncoin"
    "C34" = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
    "E34" = "  +1.34%  "
    "E35" = "  -6.07%  "
    "E36" = "  -5.38%  "
    "E37" = "  +0.06%  "
    "E38" = "  -3.39%  "
    "D39" = "0.0₃0753"
    "E39" = "  +0.01%  "
    "E40" = "  -6.50%  "
    "E41" = "  -12.02%  "
    "E42" = "  -0.25%  "
    "B43" = "Stellar"
    "C43" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    "E43" = "  -9.04%  "
    "B44" = "EnergySwap"
    "C44" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "E44" = "  -7.51%  "
    "E45" = "  +18.55%  "
    "E46" = "  -2.02%  "
    "E47" = "  -4.37%  "
    "E48" = "  -2.02%  "
    "E49" = "  -1.18%  "
    "E50" = "  -3.31%  "
    "E51" = "  -4.22%  "
}

foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}

# Cells whose new value LOOKS like a plain number (e.g. "0.609", "412.21")
# and must stay a text string (matching the source file's inlineStr cells).
# Force text entry by switching the cell to a text number format, writing
# the value, then restoring the cell style so no stray formatting sticks.
$textUpdates = @{
    "D5"  = "408.16"
    "D6"  = "132.72"
    "D8"  = "0.609"
    "D10" = "0.729"
    "D12" = "0.0000358"
    "D13" = "41.06"
    "D15" = "9.90"
    "D16" = "14.75"
    "D22" = "412.21"
    "D23" = "14.38"
    "D24" = "85.21"
    "D25" = "3.06"
    "D26" = "5.73"
    "D27" = "35.99"
    "D29" = "9.36"
    "D30" = "741.25"
    "D31" = "12.36"
    "D33" = "7.41"
    "D34" = "2.73"
    "D36" = "39.06"
    "D38" = "54.88"
    "D40" = "0.0460"
    "D42" = "0.996"
    "D43" = "0.135"
    "D44" = "27.26"
    "D45" = "3.15"
    "D46" = "145.09"
    "D47" = "3.27"
    "D48" = "2.07"
    "D49" = "2.60"
    "D51" = "2.80"
}

foreach ($addr in $textUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textUpdates[$addr]
    $cell.Style = "Normal"
}
